$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 and C1 need to hold plain-text values that look like a date and a
# time ("01/06/2012" / "00:00:50"). Assigning such strings directly to
# .Value would make Excel auto-convert them into date/time serial
# numbers. To keep them as genuine text (and avoid leaving any stray
# number-format style behind on the cell), build the text through a
# formula in a scratch cell and paste only the resulting value back.
$ws.Range("Z1").Formula = "=""01/06/2012"""
$ws.Range("Z1").Copy()
$ws.Range("B1").PasteSpecial(-4163)

$ws.Range("Z1").Formula = "=""00:00:50"""
$ws.Range("Z1").Copy()
$ws.Range("C1").PasteSpecial(-4163)

$ws.Range("Z1").Clear()

$ws.Range("D1").Value = "99aabcez"

# Remove the second and third data rows entirely, leaving only the
# header/first data row.
$ws.Rows("2:3").Delete()
